# Auto-generated edit script: update F-column (想去人数 / want-to-go count) values
# across all four worksheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 224
$ws.Range("F3").Value = 741
$ws.Range("F4").Value = 545
$ws.Range("F5").Value = 2249
$ws.Range("F6").Value = 1333
$ws.Range("F7").Value = 797
$ws.Range("F8").Value = 101
$ws.Range("F9").Value = 24
$ws.Range("F10").Value = 2901
$ws.Range("F11").Value = 25
$ws.Range("F13").Value = 1079
$ws.Range("F14").Value = 583
$ws.Range("F17").Value = 959
$ws.Range("F18").Value = 959
$ws.Range("F19").Value = 111
$ws.Range("F20").Value = 6
$ws.Range("F21").Value = 125
$ws.Range("F22").Value = 533
$ws.Range("F23").Value = 155
$ws.Range("F24").Value = 620
$ws.Range("F25").Value = 591
$ws.Range("F26").Value = 292
$ws.Range("F29").Value = 995
$ws.Range("F30").Value = 4929
$ws.Range("F31").Value = 425
$ws.Range("F32").Value = 185
$ws.Range("F33").Value = 89

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 15
$ws.Range("F5").Value = 357
$ws.Range("F6").Value = 396
$ws.Range("F11").Value = 187
$ws.Range("F14").Value = 8
$ws.Range("F18").Value = 1772
$ws.Range("F22").Value = 37
$ws.Range("F23").Value = 39
$ws.Range("F24").Value = 359
$ws.Range("F26").Value = 610
$ws.Range("F30").Value = 52
$ws.Range("F33").Value = 228
$ws.Range("F38").Value = 729
$ws.Range("F39").Value = 35

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 622
$ws.Range("F6").Value = 395
$ws.Range("F7").Value = 371

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 622
$ws.Range("F4").Value = 224
$ws.Range("F5").Value = 395
$ws.Range("F6").Value = 15
$ws.Range("F7").Value = 741
$ws.Range("F8").Value = 357
$ws.Range("F9").Value = 396
$ws.Range("F11").Value = 545
$ws.Range("F12").Value = 2249
$ws.Range("F13").Value = 1333
$ws.Range("F14").Value = 797
$ws.Range("F15").Value = 101
$ws.Range("F17").Value = 187
$ws.Range("F18").Value = 24
$ws.Range("F20").Value = 2901
$ws.Range("F21").Value = 2901
$ws.Range("F22").Value = 25
$ws.Range("F25").Value = 1079
$ws.Range("F26").Value = 583
$ws.Range("F28").Value = 371
$ws.Range("F30").Value = 959
$ws.Range("F31").Value = 959
$ws.Range("F32").Value = 111
$ws.Range("F35").Value = 125
$ws.Range("F36").Value = 155
$ws.Range("F37").Value = 37
$ws.Range("F38").Value = 39
$ws.Range("F39").Value = 620
$ws.Range("F40").Value = 591
$ws.Range("F41").Value = 359
$ws.Range("F42").Value = 610
$ws.Range("F43").Value = 292
$ws.Range("F46").Value = 995
$ws.Range("F47").Value = 4929
$ws.Range("F48").Value = 52
$ws.Range("F49").Value = 425
$ws.Range("F50").Value = 185
$ws.Range("F51").Value = 729
$ws.Range("F52").Value = 729
